$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new column at J.
#    Old layout:  I=ESTADO J=CEP K=DDD L=CELULAR M=FONE N=SITUACAO
#    New layout:  I=MUNICIPIO J=UF(new) K=CEP L=DDD M=CELULAR N=FONE O=SITUACAO
#    (old J..N shift right to K..O automatically)
# ---------------------------------------------------------------------------
$ws.Columns("J").Insert()

# ---------------------------------------------------------------------------
# 2) Header row updates
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = "MUNICIPIO"
$ws.Range("J1").Value = "UF"
$ws.Range("P1").Value = "CODIGO_MUNICIPIO"

# Give the new P1 header cell the same look (bold/border/centered) as the rest
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) MUNICIPIO (I) / UF (J) values for each data row
# ---------------------------------------------------------------------------
$ws.Range("I2").Value = "JOAO PESSOA";            $ws.Range("J2").Value = "PB"
$ws.Range("I3").Value = "CAICARA DO NORTE";        $ws.Range("J3").Value = "RN"
$ws.Range("I4").Value = "JOAO PESSOA";             $ws.Range("J4").Value = "PB"
$ws.Range("I5").Value = "BOM JESUS";                $ws.Range("J5").Value = "PB"
$ws.Range("I6").Value = "PATOS";                    $ws.Range("J6").Value = "PB"
$ws.Range("I7").Value = "BOM JESUS";                $ws.Range("J7").Value = "RN"
$ws.Range("I8").Value = "ALAGOA GRANDE";            $ws.Range("J8").Value = "PB"
$ws.Range("I9").Value = "BAYEUX";                   $ws.Range("J9").Value = "PB"
$ws.Range("I10").Value = "BELEM DO SAO FRANCISCO";  $ws.Range("J10").Value = "PE"
$ws.Range("I11").Value = "PANELAS";                 $ws.Range("J11").Value = "PE"
$ws.Range("I12").Value = "ALHANDRA";                $ws.Range("J12").Value = "PB"

# ---------------------------------------------------------------------------
# 4) Row 12 CEP is brand-new data (was blank before)
# ---------------------------------------------------------------------------
$ws.Range("K12").Value = "'56789012"

# ---------------------------------------------------------------------------
# 5) SITUACAO (O) corrections
#    rows 7,8,9 previously blank -> now CADASTRADO
#    row 12 CADASTRADO -> CADASTRANDO
# ---------------------------------------------------------------------------
$ws.Range("O7").Value = "CADASTRADO"
$ws.Range("O8").Value = "CADASTRADO"
$ws.Range("O9").Value = "CADASTRADO"
$ws.Range("O12").Value = "CADASTRANDO"

# ---------------------------------------------------------------------------
# 6) CODIGO_MUNICIPIO (P) values
# ---------------------------------------------------------------------------
$ws.Range("P2").Value = "'2507507"
$ws.Range("P3").Value = 2401859
$ws.Range("P4").Value = "'2507507"
$ws.Range("P5").Value = 2502201
$ws.Range("P6").Value = 2510808
$ws.Range("P7").Value = 2401701
$ws.Range("P8").Value = 2500304
$ws.Range("P9").Value = 2501807
$ws.Range("P10").Value = 2601607
$ws.Range("P11").Value = 2610202
